# Fix DateTime format :bug:
#
# The "Tableless" sheet stores DateTime values in column C, but the values
# previously only carried whole-day precision. Update C3/C4 to include a
# time-of-day fraction, apply the existing date+time number format (which
# already exists in the workbook's style table) to those cells, widen
# column C slightly so the longer formatted text fits, and leave the
# worksheet with C3 selected/active (as it would be after making these
# edits interactively).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tableless")

# C3: 2023-10-13 -> 2023-10-13 18:00 (45214 -> 45214.75)
$ws.Range("C3").NumberFormat = "d/m/yy h:mm;@"
$ws.Range("C3").Value = 45214.75

# C4: 2023-10-14 -> 2023-10-14 20:00 (45215 -> 45215.83333333333)
$ws.Range("C4").NumberFormat = "d/m/yy h:mm;@"
$ws.Range("C4").Value = 45215.83333333333

# Column C needs to be a bit wider to show the new date+time values;
# column D keeps its existing width.
$ws.Columns.Item(3).ColumnWidth = 12.5

# Leave the selection on C3, which also makes "Tableless" the active sheet.
$ws.Range("C3").Select() | Out-Null
